$p = $ppt.ActivePresentation

# --- Slide 3: "Source code of application" textbox -------------------------
# The box already reads "S" + "ource " + "code of application" (three runs).
# Merge the second and third runs into a single run reading
# "ource code of application" (matching a plain text edit made directly in
# that second/third run instead of touching the leading "S").
$srcSlide = $p.Slides.Item(3)
$srcShape = $srcSlide.Shapes.Item(44)           # "Textfeld 117"
$srcRange = $srcShape.TextFrame2.TextRange
$mergeChars = $srcRange.Characters(2, $srcRange.Length - 1)
$mergeChars.Text = "ource code of application"

# --- Slide 4: application-server rectangle -----------------------------
# Rename the "TomEE" label to "Payara Micro".
$appSlide = $p.Slides.Item(4)
$appShape = $appSlide.Shapes.Item(7)            # "Rechteck 5"
$appShape.TextFrame.TextRange.Text = "Payara Micro"
